$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode to "N" for rows 3 through 26 (row 2 stays "Y")
$ws.Range("D3:D26").Value = "N"

# Update the active selection to match the recorded state after the edit
$ws.Range("D4").Select()
